$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(182, 1).Value = "What is the upper limit on lithology types in a log?"
$ws.Cells.Item(182, 2).Value = "llama3.2:latest"
$ws.Cells.Item(182, 3).Value = "The upper limit on lithology types in a log is 450."

$ws.Cells.Item(183, 1).Value = "How many lithology types can a single log accommodate at maximum?"
$ws.Cells.Item(183, 2).Value = "llama3.2:latest"
$ws.Cells.Item(183, 3).Value = "The highest number of lithology types that can be represented in a log is 450."
